# Update Work Week and Social Spending
# -------------------------------------------------------------
# This workbook ("Morocco_GDPperCapita_TerritorialRef_1979_2012_CCode_504")
# has its "GDP per Capita" series on the "Data" sheet refreshed with newer
# source values, and six additional years (2011-2016) appended.
# The "Metadata" sheet (Description / Downloaded from / citations) is left
# untouched - its content is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The "Data" column (E) stores numbers as plain text in this workbook, so we
# format the relevant cells as Text before writing into them - this keeps the
# numeric-looking values ("685", "2590.76176449", ...) stored as strings
# instead of being auto-converted to real numbers by Excel.
$ws.Range("E2:E198").NumberFormat = "@"

# Refreshed GDP per Capita values for the existing years (rows 2-192).
# Only rows whose value actually changed are listed below.
$updatedValues = @{
    2   = "685"
    52  = "897"
    95  = "1132"
    132 = "2319"
    133 = "2324"
    134 = "2327"
    135 = "2340"
    136 = "2353"
    137 = "2364"
    138 = "2313"
    139 = "2263"
    140 = "2214"
    141 = "2165"
    142 = "2118"
    143 = "2138"
    144 = "2158"
    145 = "2179"
    146 = "2201"
    147 = "2222"
    148 = "2289"
    149 = "2359"
    150 = "2429"
    151 = "2503"
    152 = "2576"
    153 = "2654"
    154 = "2660"
    155 = "2700"
    156 = "2791"
    157 = "2919"
    158 = "3175"
    159 = "3268"
    160 = "3298"
    161 = "3382"
    162 = "3622"
    163 = "3457"
    164 = "3725"
    165 = "3606"
    166 = "3657"
    167 = "3784"
    168 = "3977"
    169 = "3790"
    170 = "4077"
    171 = "4073"
    172 = "4130"
    173 = "4353.7296889409"
    174 = "4194.44690978948"
    175 = "4097.28062769957"
    176 = "4461.06259109658"
    177 = "4156.92505145095"
    178 = "4603.65195843106"
    179 = "4469.18982356671"
    180 = "4729.18293904048"
    181 = "4720.77134761778"
    182 = "4753.63825365122"
    183 = "5042.07310981944"
    184 = "5140.8725323116"
    185 = "5386.85038008371"
    186 = "5584.15384967664"
    187 = "5708.03114095666"
    188 = "6076.95253367062"
    189 = "6228.43166689224"
    190 = "6530.97807194875"
    191 = "6741.46026999185"
    192 = "6931.55829013224"
}

foreach ($row in $updatedValues.Keys) {
    $ws.Range("E$row").Value = $updatedValues[$row]
}

# Append six new years of data (2011-2016) in rows 193-198.
$newRows = @(
    @{ Row = 193; Year = 2011; Value = "7226" },
    @{ Row = 194; Year = 2012; Value = "7365" },
    @{ Row = 195; Year = 2013; Value = "7619" },
    @{ Row = 196; Year = 2014; Value = "7733" },
    @{ Row = 197; Year = 2015; Value = "8001" },
    @{ Row = 198; Year = 2016; Value = "8039" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = 504.0
    $ws.Range("B$r").Value = "Morocco"
    $ws.Range("C$r").Value = "GDP per Capita"
    $ws.Range("D$r").Value = $entry.Year
    $ws.Range("E$r").Value = $entry.Value
}
